# Refresh the cryptos price/volume snapshot (scraped values) in-place.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = '@'
$ws.Range("D2").Value = '55.674.97'
$ws.Range("D2").Style = 'Normal'
$ws.Range("E2").Value = '  -2.58%  '
$ws.Range("D3").NumberFormat = '@'
$ws.Range("D3").Value = '2.959.49'
$ws.Range("D3").Style = 'Normal'
$ws.Range("E3").Value = '  -1.15%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '501.71'
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  +0.03%  '
$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '135.95'
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  -1.94%  '
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("E8").Value = '  -1.49%  '
$ws.Range("E9").Value = '  -2.81%  '
$ws.Range("E10").Value = '  -2.37%  '
$ws.Range("E11").Value = '  +0.45%  '
$ws.Range("D12").NumberFormat = '@'
$ws.Range("D12").Value = '3.475.55'
$ws.Range("D12").Style = 'Normal'
$ws.Range("E12").Value = '  -0.64%  '
$ws.Range("E13").Value = '  -1.64%  '
$ws.Range("D14").NumberFormat = '@'
$ws.Range("D14").Value = '25.73'
$ws.Range("D14").Style = 'Normal'
$ws.Range("E14").Value = '  -1.72%  '
$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '0.0000158'
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = '  -1.71%  '
$ws.Range("D16").NumberFormat = '@'
$ws.Range("D16").Value = '55.735.08'
$ws.Range("D16").Style = 'Normal'
$ws.Range("E16").Value = '  -2.54%  '
$ws.Range("D17").NumberFormat = '@'
$ws.Range("D17").Value = '2.979.77'
$ws.Range("D17").Style = 'Normal'
$ws.Range("E17").Value = '  -0.22%  '
$ws.Range("D18").NumberFormat = '@'
$ws.Range("D18").Value = '5.94'
$ws.Range("D18").Style = 'Normal'
$ws.Range("E18").Value = '  -2.45%  '
$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '12.77'
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = '  +0.57%  '
$ws.Range("D20").NumberFormat = '@'
$ws.Range("D20").Value = '7.89'
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = '  -0.20%  '
$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '324.80'
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = '  +0.60%  '
$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '0.999'
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '0.490'
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = '  -0.39%  '
$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '64.23'
$ws.Range("D24").Style = 'Normal'
$ws.Range("E24").Value = '  +0.57%  '
$ws.Range("D25").NumberFormat = '@'
$ws.Range("D25").Value = '3.094.11'
$ws.Range("D25").Style = 'Normal'
$ws.Range("E25").Value = '  -0.47%  '
$ws.Range("E26").Value = '  +0.68%  '
$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '0.161'
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  -2.09%  '
$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '0.0₃0897'
$ws.Range("D28").Style = 'Normal'
$ws.Range("E28").Value = '  -0.57%  '
$ws.Range("D29").NumberFormat = '@'
$ws.Range("D29").Value = '6.30'
$ws.Range("D29").Style = 'Normal'
$ws.Range("E29").Value = '  -4.47%  '
$ws.Range("D30").NumberFormat = '@'
$ws.Range("D30").Value = '6.91'
$ws.Range("D30").Style = 'Normal'
$ws.Range("E30").Value = '  -2.90%  '
$ws.Range("D31").NumberFormat = '@'
$ws.Range("D31").Value = '1.76'
$ws.Range("D31").Style = 'Normal'
$ws.Range("E31").Value = '  -0.73%  '
$ws.Range("D32").NumberFormat = '@'
$ws.Range("D32").Value = '20.06'
$ws.Range("D32").Style = 'Normal'
$ws.Range("E32").Value = '  -1.31%  '
$ws.Range("E33").Value = '  -3.05%  '
$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = '152.70'
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = '  -1.89%  '
$ws.Range("D35").NumberFormat = '@'
$ws.Range("D35").Value = '4.45'
$ws.Range("D35").Style = 'Normal'
$ws.Range("E35").Value = '  -3.30%  '
$ws.Range("E36").Value = '  -2.24%  '
$ws.Range("D37").NumberFormat = '@'
$ws.Range("D37").Value = '25.33'
$ws.Range("D37").Style = 'Normal'
$ws.Range("E37").Value = '  +3.42%  '
$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '1.22'
$ws.Range("D38").Style = 'Normal'
$ws.Range("E38").Value = '  -2.25%  '
$ws.Range("D39").NumberFormat = '@'
$ws.Range("D39").Value = '0.0654'
$ws.Range("D39").Style = 'Normal'
$ws.Range("E39").Value = '  -2.15%  '
$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '2.996.50'
$ws.Range("D40").Style = 'Normal'
$ws.Range("E40").Value = '  -0.87%  '
$ws.Range("D41").NumberFormat = '@'
$ws.Range("D41").Value = '36.67'
$ws.Range("D41").Style = 'Normal'
$ws.Range("E41").Value = '  -3.03%  '
$ws.Range("E42").Value = '  +0.06%  '
$ws.Range("E43").Value = '  -0.54%  '
$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '0.646'
$ws.Range("D44").Style = 'Normal'
$ws.Range("E44").Value = '  +0.41%  '
$ws.Range("D45").NumberFormat = '@'
$ws.Range("D45").Value = '2.154.08'
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = '  -2.26%  '
$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '1.33'
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  -4.16%  '
$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '5.78'
$ws.Range("D47").Style = 'Normal'
$ws.Range("E47").Value = '  -3.35%  '
$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '0.914'
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  -2.94%  '
$ws.Range("E49").Value = '  -0.90%  '
$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '19.38'
$ws.Range("D50").Style = 'Normal'
$ws.Range("E50").Value = '  -0.23%  '
$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '0.0844'
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = '  -4.07%  '
